# SystemTestCases.xlsx — "updated OPD billing test script"
#
# Content change: every data row's Status column ("Norun") is corrected to
# the properly-cased "NoRun" (rows 2-58 of Sheet1, column B).
#
# View change: the last active selection when the workbook was saved moved
# from A16 to A8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the Status column text for every test-case row (B2:B58) in one shot.
$ws.Range("B2:B58").Value = "NoRun"

# Restore the cursor/selection to A8, matching the saved view state.
$ws.Range("A8").Select()
